$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 5: targeted odds tweaks
$ws.Cells.Item(5, 7).Value = 2.6
$ws.Cells.Item(5, 9).Value = 3
$ws.Cells.Item(5, 12).Value = 3.75
$ws.Cells.Item(5, 41).Value = 17
$ws.Cells.Item(5, 42).Value = 34

# Row 6: full refresh (shift-up + per-row adjustments)
$ws.Cells.Item(6, 1).Value = 'W42VXeOQ'
$ws.Cells.Item(6, 2).Value = '11/11/2024'
$ws.Cells.Item(6, 3).Value = '21:00'
$ws.Cells.Item(6, 4).Value = 'BRAZIL - SERIE B'
$ws.Cells.Item(6, 5).Value = 'Coritiba'
$ws.Cells.Item(6, 6).Value = 'Santos'
$ws.Cells.Item(6, 7).Value = 3.4
$ws.Cells.Item(6, 8).Value = 2.9
$ws.Cells.Item(6, 9).Value = 2.38
$ws.Cells.Item(6, 10).Value = 4
$ws.Cells.Item(6, 11).Value = 1.91
$ws.Cells.Item(6, 12).Value = 3.2
$ws.Cells.Item(6, 13).Value = 1.13
$ws.Cells.Item(6, 14).Value = 6
$ws.Cells.Item(6, 15).Value = 1.5
$ws.Cells.Item(6, 16).Value = 2.5
$ws.Cells.Item(6, 17).Value = 2.6
$ws.Cells.Item(6, 18).Value = 1.48
$ws.Cells.Item(6, 19).Value = 1.57
$ws.Cells.Item(6, 20).Value = 2.25
$ws.Cells.Item(6, 21).Value = 2.1
$ws.Cells.Item(6, 22).Value = 1.67
$ws.Cells.Item(6, 23).Value = 8
$ws.Cells.Item(6, 24).Value = 15
$ws.Cells.Item(6, 25).Value = 13
$ws.Cells.Item(6, 26).Value = 34
$ws.Cells.Item(6, 27).Value = 34
$ws.Cells.Item(6, 28).Value = 41
$ws.Cells.Item(6, 29).Value = 6
$ws.Cells.Item(6, 30).Value = 5.5
$ws.Cells.Item(6, 31).Value = 17
$ws.Cells.Item(6, 32).Value = 67
$ws.Cells.Item(6, 33).Value = 6
$ws.Cells.Item(6, 34).Value = 10
$ws.Cells.Item(6, 35).Value = 10
$ws.Cells.Item(6, 36).Value = 23
$ws.Cells.Item(6, 37).Value = 23
$ws.Cells.Item(6, 38).Value = 41
$ws.Cells.Item(6, 39).Value = 351
$ws.Cells.Item(6, 40).Value = 5
$ws.Cells.Item(6, 41).Value = 21
$ws.Cells.Item(6, 42).Value = 34
$ws.Cells.Item(6, 43).Value = 67
$ws.Cells.Item(6, 44).Value = 101
$ws.Cells.Item(6, 45).Value = 351
$ws.Cells.Item(6, 46).Value = 2.25
$ws.Cells.Item(6, 47).Value = 9
$ws.Cells.Item(6, 48).Value = 67
$ws.Cells.Item(6, 49).Value = 4.33
$ws.Cells.Item(6, 50).Value = 15
$ws.Cells.Item(6, 51).Value = 29
$ws.Cells.Item(6, 52).Value = 51
$ws.Cells.Item(6, 53).Value = 81
$ws.Cells.Item(6, 54).Value = 301
$ws.Cells.Item(6, 55).Value = 81
$ws.Cells.Item(6, 56).Value = 81

# Row 7: full refresh (shift-up + per-row adjustments)
$ws.Cells.Item(7, 1).Value = 'YDlcTFWs'
$ws.Cells.Item(7, 2).Value = '11/11/2024'
$ws.Cells.Item(7, 3).Value = '21:00'
$ws.Cells.Item(7, 4).Value = 'BRAZIL - SERIE B'
$ws.Cells.Item(7, 5).Value = 'Paysandu PA'
$ws.Cells.Item(7, 6).Value = 'Brusque'
$ws.Cells.Item(7, 7).Value = 1.62
$ws.Cells.Item(7, 8).Value = 3.75
$ws.Cells.Item(7, 9).Value = 5.5
$ws.Cells.Item(7, 10).Value = 2.3
$ws.Cells.Item(7, 11).Value = 2.1
$ws.Cells.Item(7, 12).Value = 6
$ws.Cells.Item(7, 13).Value = 1.06
$ws.Cells.Item(7, 14).Value = 10
$ws.Cells.Item(7, 15).Value = 1.36
$ws.Cells.Item(7, 16).Value = 3
$ws.Cells.Item(7, 17).Value = 2.1
$ws.Cells.Item(7, 18).Value = 1.7
$ws.Cells.Item(7, 19).Value = 1.44
$ws.Cells.Item(7, 20).Value = 2.63
$ws.Cells.Item(7, 21).Value = 2.1
$ws.Cells.Item(7, 22).Value = 1.67
$ws.Cells.Item(7, 23).Value = 6
$ws.Cells.Item(7, 24).Value = 7
$ws.Cells.Item(7, 25).Value = 9
$ws.Cells.Item(7, 26).Value = 12
$ws.Cells.Item(7, 27).Value = 15
$ws.Cells.Item(7, 28).Value = 34
$ws.Cells.Item(7, 29).Value = 8.5
$ws.Cells.Item(7, 30).Value = 7
$ws.Cells.Item(7, 31).Value = 21
$ws.Cells.Item(7, 32).Value = 67
$ws.Cells.Item(7, 33).Value = 12
$ws.Cells.Item(7, 34).Value = 26
$ws.Cells.Item(7, 35).Value = 19
$ws.Cells.Item(7, 36).Value = 67
$ws.Cells.Item(7, 37).Value = 41
$ws.Cells.Item(7, 38).Value = 51
$ws.Cells.Item(7, 39).Value = 1250
$ws.Cells.Item(7, 40).Value = 3.4
$ws.Cells.Item(7, 41).Value = 8.5
$ws.Cells.Item(7, 42).Value = 23
$ws.Cells.Item(7, 43).Value = 29
$ws.Cells.Item(7, 44).Value = 51
$ws.Cells.Item(7, 45).Value = 201
$ws.Cells.Item(7, 46).Value = 2.63
$ws.Cells.Item(7, 47).Value = 9.5
$ws.Cells.Item(7, 48).Value = 67
$ws.Cells.Item(7, 49).Value = 7
$ws.Cells.Item(7, 50).Value = 34
$ws.Cells.Item(7, 51).Value = 41
$ws.Cells.Item(7, 52).Value = 126
$ws.Cells.Item(7, 53).Value = 151
$ws.Cells.Item(7, 54).Value = 351
$ws.Cells.Item(7, 55).Value = 81
$ws.Cells.Item(7, 56).Value = 81

# Row 8: full refresh (shift-up + per-row adjustments)
$ws.Cells.Item(8, 1).Value = 'lK95F3W0'
$ws.Cells.Item(8, 2).Value = '11/11/2024'
$ws.Cells.Item(8, 3).Value = '20:20'
$ws.Cells.Item(8, 4).Value = 'COLOMBIA - PRIMERA A'
$ws.Cells.Item(8, 5).Value = 'Once Caldas'
$ws.Cells.Item(8, 6).Value = 'Junior'
$ws.Cells.Item(8, 7).Value = 2.2
$ws.Cells.Item(8, 8).Value = 3
$ws.Cells.Item(8, 9).Value = 3.6
$ws.Cells.Item(8, 10).Value = 3
$ws.Cells.Item(8, 11).Value = 1.95
$ws.Cells.Item(8, 12).Value = 4.33
$ws.Cells.Item(8, 13).Value = 1.1
$ws.Cells.Item(8, 14).Value = 7
$ws.Cells.Item(8, 15).Value = 1.44
$ws.Cells.Item(8, 16).Value = 2.63
$ws.Cells.Item(8, 17).Value = 2.5
$ws.Cells.Item(8, 18).Value = 1.5
$ws.Cells.Item(8, 19).Value = 1.57
$ws.Cells.Item(8, 20).Value = 2.25
$ws.Cells.Item(8, 21).Value = 2.1
$ws.Cells.Item(8, 22).Value = 1.67
$ws.Cells.Item(8, 23).Value = 6
$ws.Cells.Item(8, 24).Value = 9.5
$ws.Cells.Item(8, 25).Value = 10
$ws.Cells.Item(8, 26).Value = 21
$ws.Cells.Item(8, 27).Value = 21
$ws.Cells.Item(8, 28).Value = 41
$ws.Cells.Item(8, 29).Value = 6.5
$ws.Cells.Item(8, 30).Value = 6
$ws.Cells.Item(8, 31).Value = 19
$ws.Cells.Item(8, 32).Value = 67
$ws.Cells.Item(8, 33).Value = 8.5
$ws.Cells.Item(8, 34).Value = 17
$ws.Cells.Item(8, 35).Value = 13
$ws.Cells.Item(8, 36).Value = 41
$ws.Cells.Item(8, 37).Value = 34
$ws.Cells.Item(8, 38).Value = 41
$ws.Cells.Item(8, 39).Value = 900
$ws.Cells.Item(8, 40).Value = 4
$ws.Cells.Item(8, 41).Value = 13
$ws.Cells.Item(8, 42).Value = 29
$ws.Cells.Item(8, 43).Value = 41
$ws.Cells.Item(8, 44).Value = 81
$ws.Cells.Item(8, 45).Value = 251
$ws.Cells.Item(8, 46).Value = 2.25
$ws.Cells.Item(8, 47).Value = 9
$ws.Cells.Item(8, 48).Value = 67
$ws.Cells.Item(8, 49).Value = 5
$ws.Cells.Item(8, 50).Value = 21
$ws.Cells.Item(8, 51).Value = 34
$ws.Cells.Item(8, 52).Value = 81
$ws.Cells.Item(8, 53).Value = 126
$ws.Cells.Item(8, 54).Value = 351
$ws.Cells.Item(8, 55).Value = 126
$ws.Cells.Item(8, 56).Value = 126

# Row 9: full refresh (shift-up + per-row adjustments)
$ws.Cells.Item(9, 1).Value = 'ARJPKb8t'
$ws.Cells.Item(9, 2).Value = '11/11/2024'
$ws.Cells.Item(9, 3).Value = '22:30'
$ws.Cells.Item(9, 4).Value = 'COLOMBIA - PRIMERA A'
$ws.Cells.Item(9, 5).Value = 'America De Cali'
$ws.Cells.Item(9, 6).Value = 'Santa Fe'
$ws.Cells.Item(9, 7).Value = 1.83
$ws.Cells.Item(9, 8).Value = 3.1
$ws.Cells.Item(9, 9).Value = 5
$ws.Cells.Item(9, 10).Value = 2.6
$ws.Cells.Item(9, 11).Value = 2
$ws.Cells.Item(9, 12).Value = 5.5
$ws.Cells.Item(9, 13).Value = 1.11
$ws.Cells.Item(9, 14).Value = 6.5
$ws.Cells.Item(9, 15).Value = 1.44
$ws.Cells.Item(9, 16).Value = 2.63
$ws.Cells.Item(9, 17).Value = 2.4
$ws.Cells.Item(9, 18).Value = 1.53
$ws.Cells.Item(9, 19).Value = 1.53
$ws.Cells.Item(9, 20).Value = 2.38
$ws.Cells.Item(9, 21).Value = 2.2
$ws.Cells.Item(9, 22).Value = 1.62
$ws.Cells.Item(9, 23).Value = 5.5
$ws.Cells.Item(9, 24).Value = 7.5
$ws.Cells.Item(9, 25).Value = 9.5
$ws.Cells.Item(9, 26).Value = 15
$ws.Cells.Item(9, 27).Value = 19
$ws.Cells.Item(9, 28).Value = 41
$ws.Cells.Item(9, 29).Value = 6.5
$ws.Cells.Item(9, 30).Value = 6.5
$ws.Cells.Item(9, 31).Value = 19
$ws.Cells.Item(9, 32).Value = 81
$ws.Cells.Item(9, 33).Value = 10
$ws.Cells.Item(9, 34).Value = 23
$ws.Cells.Item(9, 35).Value = 17
$ws.Cells.Item(9, 36).Value = 51
$ws.Cells.Item(9, 37).Value = 41
$ws.Cells.Item(9, 38).Value = 51
$ws.Cells.Item(9, 39).Value = 201
$ws.Cells.Item(9, 40).Value = 3.6
$ws.Cells.Item(9, 41).Value = 10
$ws.Cells.Item(9, 42).Value = 26
$ws.Cells.Item(9, 43).Value = 41
$ws.Cells.Item(9, 44).Value = 67
$ws.Cells.Item(9, 45).Value = 251
$ws.Cells.Item(9, 46).Value = 2.38
$ws.Cells.Item(9, 47).Value = 9.5
$ws.Cells.Item(9, 48).Value = 81
$ws.Cells.Item(9, 49).Value = 6
$ws.Cells.Item(9, 50).Value = 29
$ws.Cells.Item(9, 51).Value = 41
$ws.Cells.Item(9, 52).Value = 101
$ws.Cells.Item(9, 53).Value = 151
$ws.Cells.Item(9, 54).Value = 351
$ws.Cells.Item(9, 55).Value = 126
$ws.Cells.Item(9, 56).Value = 126

# Row 10: full refresh (shift-up + per-row adjustments)
$ws.Cells.Item(10, 1).Value = 'SlBBcd9a'
$ws.Cells.Item(10, 2).Value = '11/11/2024'
$ws.Cells.Item(10, 3).Value = '20:00'
$ws.Cells.Item(10, 4).Value = 'PARAGUAY - PRIMERA DIVISION'
$ws.Cells.Item(10, 5).Value = 'Ameliano'
$ws.Cells.Item(10, 6).Value = '2 de Mayo'
$ws.Cells.Item(10, 7).Value = 2.38
$ws.Cells.Item(10, 8).Value = 2.8
$ws.Cells.Item(10, 9).Value = 3.25
$ws.Cells.Item(10, 10).Value = 3.4
$ws.Cells.Item(10, 11).Value = 1.83
$ws.Cells.Item(10, 12).Value = 4.33
$ws.Cells.Item(10, 13).Value = 1.13
$ws.Cells.Item(10, 14).Value = 6
$ws.Cells.Item(10, 15).Value = 1.57
$ws.Cells.Item(10, 16).Value = 2.25
$ws.Cells.Item(10, 17).Value = 2.88
$ws.Cells.Item(10, 18).Value = 1.4
$ws.Cells.Item(10, 19).Value = 1.62
$ws.Cells.Item(10, 20).Value = 2.2
$ws.Cells.Item(10, 21).Value = 2.25
$ws.Cells.Item(10, 22).Value = 1.57
$ws.Cells.Item(10, 23).Value = 6
$ws.Cells.Item(10, 24).Value = 10
$ws.Cells.Item(10, 25).Value = 11
$ws.Cells.Item(10, 26).Value = 23
$ws.Cells.Item(10, 27).Value = 26
$ws.Cells.Item(10, 28).Value = 41
$ws.Cells.Item(10, 29).Value = 5.5
$ws.Cells.Item(10, 30).Value = 6
$ws.Cells.Item(10, 31).Value = 21
$ws.Cells.Item(10, 32).Value = 81
$ws.Cells.Item(10, 33).Value = 7
$ws.Cells.Item(10, 34).Value = 15
$ws.Cells.Item(10, 35).Value = 13
$ws.Cells.Item(10, 36).Value = 34
$ws.Cells.Item(10, 37).Value = 34
$ws.Cells.Item(10, 38).Value = 51
$ws.Cells.Item(10, 39).Value = 1250
$ws.Cells.Item(10, 40).Value = 4.33
$ws.Cells.Item(10, 41).Value = 15
$ws.Cells.Item(10, 42).Value = 34
$ws.Cells.Item(10, 43).Value = 51
$ws.Cells.Item(10, 44).Value = 101
$ws.Cells.Item(10, 45).Value = 351
$ws.Cells.Item(10, 46).Value = 2.2
$ws.Cells.Item(10, 47).Value = 9.5
$ws.Cells.Item(10, 48).Value = 81
$ws.Cells.Item(10, 49).Value = 5
$ws.Cells.Item(10, 50).Value = 21
$ws.Cells.Item(10, 51).Value = 41
$ws.Cells.Item(10, 52).Value = 81
$ws.Cells.Item(10, 53).Value = 126
$ws.Cells.Item(10, 54).Value = 351
$ws.Cells.Item(10, 55).Value = 51
$ws.Cells.Item(10, 56).Value = 51

# Remove the now-obsolete last row (data shifted up by one)
$ws.Rows.Item(11).Delete()
